# Update Digikey BOM comments
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the BOM header/date comment in C1
$ws.Range("C1").Value = "Uzebox Omega v1.1.1 Digi-Key BOM. Last updated 3rd September 2025."

# Update part description comments
$ws.Range("C3").Value = "CONN POWER JACK 2.1MM (7-12VDC)"
$ws.Range("C7").Value = "SOCKET IC OPEN FRAME 40POS .6"" (MCU)"
$ws.Range("C11").Value = "CAP ALUM 1UF 20% 350V RADIAL TH         (C16,C20)"

# Move the active cell selection to A11
$ws.Range("A11").Select()
